# Word COM-interop script implementing the diff:
#  - move the "_GoBack" bookmark from the "Laser Cat" bullet down to the
#    "Gameplay gets faster after a set interval" bullet in the
#    "Rating of Task Difficulty" section
#  - rewrite several bullet texts (lives -> health, randomized timings ->
#    timing and areas, progressive difficulty wording, etc.)
#  - drop the two now-stale <w:lastRenderedPageBreak/> markers

$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $rng = $d.Paragraphs($paraIndex).Range
    $ok = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $newText, 2)
    return $ok
}

# ---------------------------------------------------------------------
# 1. Remove the _GoBack bookmark currently sitting after "Laser Cat".
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2. Plain text rewrites (paragraph-scoped so duplicate phrases such as
#    "Lives and scores" elsewhere in the document are left untouched).
# ---------------------------------------------------------------------

# "Player has up to 5 lives ..." bullet under Game Description
Replace-InParagraph 6 `
    "Player has up to 5 lives and the game ends when all lives are gone." `
    "Player has up to a set number of health and the game ends when health is completely depleted."

# "The more cycles are completed ..." bullet under Features of Game
Replace-InParagraph 16 `
    "The more cycles are completed, the faster the game gets (progressive difficulty)." `
    "As time goes by, the cat paws should move at a faster rate"

# ---- List of Tasks to Achieve Features ----
Replace-InParagraph 21 "Lives and scores" "Health and scores"

Replace-InParagraph 24 `
    "Create a dot to represent the player (black dot if possible)" `
    "Create a dot (or a set of dots) to represent the player"

Replace-InParagraph 28 `
    "Make the cat paws appear at randomized times" `
    "Make the cat paws appear at randomized times and areas"

Replace-InParagraph 29 "Lives and scores" "Health and scores"

Replace-InParagraph 30 `
    "Game ends when player loses all lives" `
    "Game ends when player loses all health points"

Replace-InParagraph 31 `
    "Gameplay gets faster after every interval of 30 seconds" `
    "Gameplay gets faster at a set interval"

# ---- Rating of Task Difficulty heading: drop stray page-break marker ----
Replace-InParagraph 39 "Rating of Task Difficulty" "Rating of Task Difficulty"

# ---- Rating of Task Difficulty > Main Features ----
Replace-InParagraph 47 "Lives and scores" "Health and scores"

Replace-InParagraph 51 `
    "Create a dot to represent the player (black dot if possible)" `
    "Create a dot (or a set of dots) to represent the player"

Replace-InParagraph 58 `
    "Make the cat paw appear at randomized timings" `
    "Make the cat paw appear at randomized timing and areas"

Replace-InParagraph 60 "Lives and scores" "Health and scores"

Replace-InParagraph 61 `
    "Game ends when player loses all lives" `
    "Game ends when player loses all health points"

Replace-InParagraph 63 `
    "Gameplay gets faster after 30 seconds intervals" `
    "Gameplay gets faster after a set interval"

# ---------------------------------------------------------------------
# 3. Re-add the _GoBack bookmark right after the text we just placed in
#    paragraph 63 ("Gameplay gets faster after a set interval").
#    A temporary marker character is used so the zero-length bookmark
#    lands after the run rather than before it.
# ---------------------------------------------------------------------
$tailRng = $d.Content
$tailRng.Find.Execute("Gameplay gets faster after a set interval")
$tailRng.Collapse(0)
$tailRng.InsertAfter("@@MARK@@")

$markRng = $d.Content
$markRng.Find.Execute("@@MARK@@")
$markRng.Collapse(1)
$d.Bookmarks.Add("_GoBack", $markRng)

$markRng2 = $d.Content
$markRng2.Find.Execute("@@MARK@@")
$markRng2.Text = ""

# ---- Estimation of Time Required to Complete Tasks heading: drop page-break ----
Replace-InParagraph 79 `
    "Estimation of Time Required to Complete Tasks" `
    "Estimation of Time Required to Complete Tasks"

# ---- Estimation of Time Required to Complete Tasks > Main Features ----
Replace-InParagraph 91 "Lives and scores" "Health and scores"

Replace-InParagraph 95 `
    "Create a dot to represent the player (black dot if possible)" `
    "Create a dot (or a set of dots) to represent the player"

Replace-InParagraph 102 `
    "Make the cat paw appear at randomized timings" `
    "Make the cat paw appear at randomized timing and areas"

Replace-InParagraph 104 "Lives and scores" "Health and scores"

Replace-InParagraph 105 `
    "Game ends when player loses all lives" `
    "Game ends when player loses all health points"

Replace-InParagraph 107 `
    "Gameplay gets faster after 30 seconds intervals" `
    "Gameplay gets faster after a set interval"

Write-Output "done"
